$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '29.350.18'
Set-TextValue 'E2' '  -0.20%  '
Set-TextValue 'D3' '1.843.34'
Set-TextValue 'E3' '  -0.32%  '
Set-TextValue 'E4' '  -0.27%  '
Set-TextValue 'D5' '240.54'
Set-TextValue 'E5' '  -0.02%  '
Set-TextValue 'D6' '0.6321'
Set-TextValue 'E6' '  +0.59%  '
Set-TextValue 'D7' '0.9989'
Set-TextValue 'E7' '  -0.19%  '
Set-TextValue 'D8' '0.07478'
Set-TextValue 'E8' '  -2.33%  '
Set-TextValue 'D9' '0.2905'
Set-TextValue 'E9' '  -0.19%  '
Set-TextValue 'D10' '24.38'
Set-TextValue 'E10' '  -1.79%  '
Set-TextValue 'D11' '0.07724'
Set-TextValue 'E11' '  -0.18%  '
Set-TextValue 'D12' '1.843.74'
Set-TextValue 'E12' '  -2.36%  '
Set-TextValue 'D13' '5.001'
Set-TextValue 'E13' '  -0.83%  '
Set-TextValue 'D14' '0.6778'
Set-TextValue 'E14' '  -0.45%  '
Set-TextValue 'D15' '0.00001021'
Set-TextValue 'E15' '  -5.09%  '
Set-TextValue 'D16' '82.17'
Set-TextValue 'E16' '  -1.41%  '
Set-TextValue 'D17' '6.137'
Set-TextValue 'E17' '  -0.64%  '
Set-TextValue 'D18' '29.389.67'
Set-TextValue 'E18' '  -0.09%  '
Set-TextValue 'D19' '228.59'
Set-TextValue 'E20' '  -0.24%  '
Set-TextValue 'D21' '0.9987'
Set-TextValue 'E21' '  -0.23%  '
Set-TextValue 'D22' '7.435'
Set-TextValue 'E22' '  -0.58%  '
Set-TextValue 'D23' '0.9983'
Set-TextValue 'E23' '  -0.28%  '
Set-TextValue 'D24' '158.73'
Set-TextValue 'E24' '  +0.42%  '
Set-TextValue 'D25' '0.1379'
Set-TextValue 'E25' '  -0.24%  '
Set-TextValue 'D26' '8.410'
Set-TextValue 'E26' '  -0.13%  '
Set-TextValue 'D27' '17.56'
Set-TextValue 'E27' '  -0.85%  '
Set-TextValue 'D28' '0.06320'
Set-TextValue 'E28' '  +12.94%  '
Set-TextValue 'D29' '1.379'
Set-TextValue 'E29' '  -0.42%  '
Set-TextValue 'D30' '1.474'
Set-TextValue 'E30' '  +0.79%  '
Set-TextValue 'D31' '4.089'
Set-TextValue 'E31' '  -1.07%  '
Set-TextValue 'D32' '4.049'
Set-TextValue 'E32' '  -0.43%  '
Set-TextValue 'D33' '1.819'
Set-TextValue 'E33' '  -1.25%  '
Set-TextValue 'D34' '1.140'
Set-TextValue 'E34' '  -2.11%  '
Set-TextValue 'D35' '0.6976'
Set-TextValue 'E35' '  +0.15%  '
Set-TextValue 'D36' '2.577'
Set-TextValue 'E36' '  -0.42%  '
Set-TextValue 'D37' '2.833'
Set-TextValue 'E37' '  +4.24%  '
Set-TextValue 'D38' '1.253.25'
Set-TextValue 'E38' '  +1.92%  '
Set-TextValue 'D39' '0.01819'
Set-TextValue 'E39' '  +0.81%  '
Set-TextValue 'D40' '6.558'
Set-TextValue 'E40' '  +1.94%  '
Set-TextValue 'D41' '0.9072'
Set-TextValue 'E41' '  -0.02%  '
Set-TextValue 'D42' '0.9980'
Set-TextValue 'E42' '  -0.29%  '
Set-TextValue 'D43' '2.007.69'
Set-TextValue 'E43' '  -18.41%  '
Set-TextValue 'D44' '101.31'
Set-TextValue 'E44' '  -0.48%  '
Set-TextValue 'D45' '66.27'
Set-TextValue 'E45' '  +0.35%  '
Set-TextValue 'E46' '  -0.56%  '
Set-TextValue 'D47' '0.1174'
Set-TextValue 'E47' '  +2.08%  '
Set-TextValue 'D48' '7.044'
Set-TextValue 'E48' '  -2.18%  '
Set-TextValue 'D49' '9.035'
Set-TextValue 'E49' '  +0.07%  '
Set-TextValue 'D50' '1.686'
Set-TextValue 'E50' '  +0.06%  '
Set-TextValue 'D51' '0.3930'
Set-TextValue 'E51' '  -2.17%  '
